$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# The "Definition" column (D) was left empty for the concept rows; fill it
# with the same text as the "Display" column (C) for each concept (rows 2-6).
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($row, 3).Value2
}
